$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the time-slot values in column C (rows 8 and 9) to the new values
$ws.Range("C8").Value = "13:5-13:10"
$ws.Range("C9").Value = "13:10-13:15"

# Move the active selection from C16 to C18
$ws.Range("C18").Select()
